$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.062.48"
$ws.Range("E2").Value = "  -6.22%  "
$ws.Range("D3").Value = "3.331.04"
$ws.Range("E3").Value = "  -8.09%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'181.68"
$ws.Range("E5").Value = "  -9.77%  "
$ws.Range("D6").Value = "'518.30"
$ws.Range("E6").Value = "  -10.47%  "
$ws.Range("D7").Value = "'0.586"
$ws.Range("E7").Value = "  -5.77%  "
$ws.Range("D8").Value = "3.321.75"
$ws.Range("E8").Value = "  -8.29%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'0.607"
$ws.Range("E10").Value = "  -11.33%  "
$ws.Range("D11").Value = "'55.96"
$ws.Range("E11").Value = "  -7.79%  "
$ws.Range("D12").Value = "'0.128"
$ws.Range("E12").Value = "  -14.22%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = "  -13.72%  "
$ws.Range("D14").Value = "'9.00"
$ws.Range("E14").Value = "  -11.58%  "
$ws.Range("D15").Value = "3.888.85"
$ws.Range("E15").Value = "  -7.44%  "
$ws.Range("E16").Value = "  -4.71%  "
$ws.Range("D17").Value = "3.346.71"
$ws.Range("E17").Value = "  -7.79%  "
$ws.Range("D18").Value = "63.862.68"
$ws.Range("E18").Value = "  -6.41%  "
$ws.Range("D19").Value = "'16.94"
$ws.Range("E19").Value = "  -12.39%  "
$ws.Range("D20").Value = "'10.74"
$ws.Range("E20").Value = "  -13.69%  "
$ws.Range("D21").Value = "'0.941"
$ws.Range("E21").Value = "  -12.58%  "
$ws.Range("D22").Value = "'364.41"
$ws.Range("E22").Value = "  -10.15%  "
$ws.Range("D23").Value = "'79.12"
$ws.Range("E23").Value = "  -7.66%  "
$ws.Range("D24").Value = "'3.62"
$ws.Range("E24").Value = "  -14.91%  "
$ws.Range("D25").Value = "'10.45"
$ws.Range("E25").Value = "  -19.47%  "
$ws.Range("D26").Value = "'5.91"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("D27").Value = "'3.64"
$ws.Range("E27").Value = "  -9.07%  "
$ws.Range("D28").Value = "'2.57"
$ws.Range("E28").Value = "  -12.30%  "
$ws.Range("D29").Value = "'10.98"
$ws.Range("E29").Value = "  -13.04%  "
$ws.Range("D30").Value = "'8.21"
$ws.Range("E30").Value = "  -12.71%  "
$ws.Range("D31").Value = "'650.90"
$ws.Range("E31").Value = "  -4.67%  "
$ws.Range("D32").Value = "'28.30"
$ws.Range("E32").Value = "  -11.24%  "
$ws.Range("D33").Value = "'6.52"
$ws.Range("E33").Value = "  -16.53%  "
$ws.Range("D34").Value = "'10.85"
$ws.Range("E34").Value = "  -11.52%  "
$ws.Range("D35").Value = "'58.03"
$ws.Range("E35").Value = "  -9.17%  "
$ws.Range("D36").Value = "'0.101"
$ws.Range("E36").Value = "  -11.70%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'35.26"
$ws.Range("E38").Value = "  -15.81%  "
$ws.Range("D39").Value = "'0.366"
$ws.Range("E39").Value = "  -12.20%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'0.122"
$ws.Range("E41").Value = "  -10.09%  "
$ws.Range("D42").Value = "'2.67"
$ws.Range("E42").Value = "  -16.17%  "
$ws.Range("D43").Value = "2.736.42"
$ws.Range("E43").Value = "  -14.49%  "
$ws.Range("D44").Value = "'2.56"
$ws.Range("E44").Value = "  -9.66%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0597"
$ws.Range("E45").Value = "  -22.28%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0376"
$ws.Range("E46").Value = "  -9.90%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'133.25"
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("D48").Value = "'0.122"
$ws.Range("E48").Value = "  -7.73%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").Value = "'2.23"
$ws.Range("E49").Value = "  -17.37%  "
$ws.Range("D50").Value = "'2.79"
$ws.Range("E50").Value = "  -9.93%  "
$ws.Range("D51").Value = "'2.51"
$ws.Range("E51").Value = "  -7.97%  "
